# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# 1. Bump the "Date" metadata value on the Metadata sheet.
# 2. Add a new "Mapping: Spécification métier vers l'extension ROR
#    ReopeningDate" column (AL) to the Elements sheet, with the mapping
#    value for Extension.value[x] filled in.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "Date" metadata value on the Metadata sheet ---
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- 2. Add new mapping column on the Elements sheet ---
$elemWs = $wb.Worksheets.Item("Elements")

# Header (row 1) + data rows (2-6) values for the new column AL (38)
$elemWs.Cells.Item(1, 38).Value = "Mapping: Spécification métier vers l'extension ROR ReopeningDate"
$elemWs.Cells.Item(2, 38).Value = ""
$elemWs.Cells.Item(3, 38).Value = ""
$elemWs.Cells.Item(4, 38).Value = ""
$elemWs.Cells.Item(5, 38).Value = ""
$elemWs.Cells.Item(6, 38).Value = "datePrevisionnelleReouverture"

# Match formatting of the neighbouring "Mapping: RIM Mapping" column (AK)
$elemWs.Range("AK1").Copy()
$elemWs.Range("AL1").PasteSpecial(-4122)
$elemWs.Range("AK2:AK6").Copy()
$elemWs.Range("AL2:AL6").PasteSpecial(-4122)

# New column width (bestFit-style width computed for the new header/content)
$elemWs.Columns.Item(38).ColumnWidth = 70.0
